$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0056609618364461753
$ws.Range("D2").Value = 0.073217204090059246
$ws.Range("E2").Value = 0.095408322749467839
$ws.Range("C3").Value = 0.0037979955943588145
$ws.Range("D3").Value = 0.050673817972811667
$ws.Range("E3").Value = 0.065561974972113163
$ws.Range("C4").Value = 0.0070964437292743579
$ws.Range("D4").Value = 0.13859979494708694
$ws.Range("E4").Value = 0.16641804022163989
$ws.Range("C5").Value = 0.0050843459837167588
$ws.Range("D5").Value = 0.086801200059406933
$ws.Range("E5").Value = 0.106731855417927
$ws.Range("C6").Value = 0.007542896022194252
$ws.Range("D6").Value = 0.19496564890623913
$ws.Range("E6").Value = 0.22453399886161962
$ws.Range("C7").Value = 0.005333893564933097
$ws.Range("D7").Value = 0.12610732505419078
$ws.Range("E7").Value = 0.14701620786865185
$ws.Range("C8").Value = 0.0079659267190836745
$ws.Range("D8").Value = 0.2329418081132954
$ws.Range("E8").Value = 0.26416844947965257
$ws.Range("C9").Value = 0.0065003685768248278
$ws.Range("D9").Value = 0.16654005299830493
$ws.Range("E9").Value = 0.19202152224193431
$ws.Range("C10").Value = 0.0093611833809891555
$ws.Range("D10").Value = 0.28233688043287641
$ws.Range("E10").Value = 0.31903296445566243
$ws.Range("C11").Value = 0.0075037762274986989
$ws.Range("D11").Value = 0.18359997360452074
$ws.Range("E11").Value = 0.21301480460868555
$ws.Range("C12").Value = 0.01014152948556128
$ws.Range("D12").Value = 0.29328395431064341
$ws.Range("E12").Value = 0.33303901550060971
$ws.Range("C13").Value = 0.0084919689494189713
$ws.Range("D13").Value = 0.21551563373440893
$ws.Range("E13").Value = 0.24880418392123113
$ws.Range("C14").Value = 0.010856566510857141
$ws.Range("D14").Value = 0.28975888315172121
$ws.Range("E14").Value = 0.33231690820766052
$ws.Range("C15").Value = 0.0084218199647003462
$ws.Range("D15").Value = 0.22547259788089202
$ws.Range("E15").Value = 0.25848616378406114
$ws.Range("C16").Value = 0.011107220759918631
$ws.Range("D16").Value = 0.20691224241326073
$ws.Range("E16").Value = 0.25045283869015361
$ws.Range("C17").Value = 0.010200149358761684
$ws.Range("D17").Value = 0.19518954905737837
$ws.Range("E17").Value = 0.23517417286661321
$ws.Range("C18").Value = 0.010415933085649904
$ws.Range("D18").Value = 0.070437645221424219
$ws.Range("E18").Value = 0.1112683757103657
$ws.Range("C19").Value = 0.011149136666402118
$ws.Range("D19").Value = 0.12182423556226543
$ws.Range("E19").Value = 0.16552889318288241
